$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 26.25523382790317
$ws.Range("C2").Value = 8.372093544870099
$ws.Range("D2").Value = 4.671133373143253
$ws.Range("F2").Value = 66.98012672441862
$ws.Range("G2").Value = 3.817106714637992
$ws.Range("J2").Value = 10.41291891820497
$ws.Range("L2").Value = 11.3833242719874
$ws.Range("M2").Value = 20.99770374164495
$ws.Range("B3").Value = 26.17659437665522
$ws.Range("C3").Value = 8.14090682226276
$ws.Range("D3").Value = 4.508045296424766
$ws.Range("F3").Value = 66.22307749789003
$ws.Range("G3").Value = 3.821938050124239
$ws.Range("J3").Value = 10.40548065560088
$ws.Range("L3").Value = 11.42366828085592
$ws.Range("M3").Value = 21.04601027226497
$ws.Range("B4").Value = 26.13787959771213
$ws.Range("C4").Value = 8.000865743225051
$ws.Range("D4").Value = 4.404885251336276
$ws.Range("F4").Value = 65.76127142301665
$ws.Range("G4").Value = 3.825054301332016
$ws.Range("J4").Value = 10.40086469634724
$ws.Range("L4").Value = 11.45015295067449
$ws.Range("M4").Value = 21.08095969288446
$ws.Range("B5").Value = 26.12452046052153
$ws.Range("C5").Value = 7.944400363544128
$ws.Range("D5").Value = 4.362131572787106
$ws.Range("F5").Value = 65.57397218575311
$ws.Range("G5").Value = 3.82636202916249
$ws.Range("J5").Value = 10.39897091894476
$ws.Range("L5").Value = 11.46137736884005
$ws.Range("M5").Value = 21.09652917313276
$ws.Range("B6").Value = 26.12244844827935
$ws.Range("C6").Value = 7.935064351474429
$ws.Range("D6").Value = 4.354990571144456
$ws.Range("F6").Value = 65.54292880850015
$ws.Range("G6").Value = 3.82658146602233
$ws.Range("J6").Value = 10.39865567605124
$ws.Range("L6").Value = 11.46326727457559
$ws.Range("M6").Value = 21.09919455617031
$ws.Range("B7").Value = 26.13768963149429
$ws.Range("C7").Value = 8.000101627940301
$ws.Range("D7").Value = 4.404311492913549
$ws.Range("F7").Value = 65.75874166512615
$ws.Range("G7").Value = 3.825071784415504
$ws.Range("J7").Value = 10.40083920837755
$ws.Range("L7").Value = 11.45030257797431
$ws.Range("M7").Value = 21.08116429780328
$ws.Range("B8").Value = 26.22614024111531
$ws.Range("C8").Value = 8.292056867682833
$ws.Range("D8").Value = 4.615553783183287
$ws.Range("F8").Value = 66.71854689545448
$ws.Range("G8").Value = 3.818741568744097
$ws.Range("J8").Value = 10.41036406533936
$ws.Range("L8").Value = 11.39687982169001
$ws.Range("M8").Value = 21.01326058184637
$ws.Range("B9").Value = 26.47488899631536
$ws.Range("C9").Value = 8.874520544536379
$ws.Range("D9").Value = 5.004053319458113
$ws.Range("F9").Value = 68.61893146090921
$ws.Range("G9").Value = 3.807509084499478
$ws.Range("J9").Value = 10.42867583173773
$ws.Range("L9").Value = 11.30567295429397
$ws.Range("M9").Value = 20.92217611736386
$ws.Range("B10").Value = 26.70247326745446
$ws.Range("C10").Value = 9.30192017593766
$ws.Range("D10").Value = 5.271735401615669
$ws.Range("F10").Value = 70.01841300489924
$ws.Range("G10").Value = 3.799966095488647
$ws.Range("J10").Value = 10.44193588648783
$ws.Range("L10").Value = 11.24687344269521
$ws.Range("M10").Value = 20.88103898717566
$ws.Range("B11").Value = 26.81544633494963
$ws.Range("C11").Value = 9.494975631217434
$ws.Range("D11").Value = 5.389302274004
$ws.Range("F11").Value = 70.6541117059277
$ws.Range("G11").Value = 3.796686409461179
$ws.Range("J11").Value = 10.44793264035041
$ws.Range("L11").Value = 11.221895683884
$ws.Range("M11").Value = 20.86794497122417
$ws.Range("B12").Value = 26.85955493891611
$ws.Range("C12").Value = 9.567788852125776
$ws.Range("D12").Value = 5.433192509071864
$ws.Range("F12").Value = 70.89456022584753
$ws.Range("G12").Value = 3.795466111174929
$ws.Range("J12").Value = 10.45019888203991
$ws.Range("L12").Value = 11.2126910050112
$ws.Range("M12").Value = 20.86379584135911
$ws.Range("B13").Value = 26.84999679596218
$ws.Range("C13").Value = 9.552121676083067
$ws.Range("D13").Value = 5.423768341536711
$ws.Range("F13").Value = 70.84278989653447
$ws.Range("G13").Value = 3.795727964015219
$ws.Range("J13").Value = 10.44971100324435
$ws.Range("L13").Value = 11.21466211843653
$ws.Range("M13").Value = 20.86465342188553
$ws.Range("B14").Value = 26.81904873882218
$ws.Range("C14").Value = 9.500972403569049
$ws.Range("D14").Value = 5.392925905835185
$ws.Range("F14").Value = 70.67389990813852
$ws.Range("G14").Value = 3.796585581855537
$ws.Range("J14").Value = 10.44811917483536
$ws.Range("L14").Value = 11.2211333258881
$ws.Range("M14").Value = 20.86758739480001
$ws.Range("B15").Value = 26.80026419254239
$ws.Range("C15").Value = 9.46960118080082
$ws.Range("D15").Value = 5.373951308472477
$ws.Range("F15").Value = 70.57040966515444
$ws.Range("G15").Value = 3.797113711942109
$ws.Range("J15").Value = 10.4471435493103
$ws.Range("L15").Value = 11.22513016545976
$ws.Range("M15").Value = 20.86948996109857
$ws.Range("B16").Value = 26.69527775663207
$ws.Range("C16").Value = 9.289267611669947
$ws.Range("D16").Value = 5.263965307835576
$ws.Range("F16").Value = 69.97683759819768
$ws.Range("G16").Value = 3.800183468696493
$ws.Range("J16").Value = 10.44154331486716
$ws.Range("L16").Value = 11.24854136091173
$ws.Range("M16").Value = 20.88200790249292
$ws.Range("B17").Value = 26.63327041175279
$ws.Range("C17").Value = 9.178217182755082
$ws.Range("D17").Value = 5.195396832946585
$ws.Range("F17").Value = 69.61235947637782
$ws.Range("G17").Value = 3.802105393918765
$ws.Range("J17").Value = 10.43809904859986
$ws.Range("L17").Value = 11.2633562977504
$ws.Range("M17").Value = 20.89112751820399
$ws.Range("B18").Value = 26.59849629561436
$ws.Range("C18").Value = 9.114220102415221
$ws.Range("D18").Value = 5.155563747412704
$ws.Range("F18").Value = 69.40264198308942
$ws.Range("G18").Value = 3.803225118793653
$ws.Range("J18").Value = 10.4361146453289
$ws.Range("L18").Value = 11.27204414201015
$ws.Range("M18").Value = 20.89690176179224
$ws.Range("B19").Value = 26.58687624916565
$ws.Range("C19").Value = 9.092533569170371
$ws.Range("D19").Value = 5.142010053642309
$ws.Range("F19").Value = 69.33162599611418
$ws.Range("G19").Value = 3.803606696631688
$ws.Range("J19").Value = 10.43544217014154
$ws.Range("L19").Value = 11.27501434745328
$ws.Range("M19").Value = 20.89894761419858
$ws.Range("B20").Value = 26.63977920553288
$ws.Range("C20").Value = 9.190052222228305
$ws.Range("D20").Value = 5.202737068354421
$ws.Range("F20").Value = 69.65116804304635
$ws.Range("G20").Value = 3.801899324488593
$ws.Range("J20").Value = 10.43846604012029
$ws.Range("L20").Value = 11.2617619765657
$ws.Range("M20").Value = 20.89010197230859
$ws.Range("B21").Value = 26.82810313905123
$ws.Range("C21").Value = 9.516004850191916
$ws.Range("D21").Value = 5.40200235535926
$ws.Range("F21").Value = 70.72351563129972
$ws.Range("G21").Value = 3.796333092475381
$ws.Range("J21").Value = 10.44858685441439
$ws.Range("L21").Value = 11.21922569155366
$ws.Range("M21").Value = 20.86670364408194
$ws.Range("B22").Value = 26.95890950077262
$ws.Range("C22").Value = 9.727289057722478
$ws.Range("D22").Value = 5.528553312915242
$ws.Range("F22").Value = 71.42267875862291
$ws.Range("G22").Value = 3.792821342570813
$ws.Range("J22").Value = 10.45517480944017
$ws.Range("L22").Value = 11.19290511040212
$ws.Range("M22").Value = 20.85612894440244
$ws.Range("B23").Value = 26.88839932891142
$ws.Range("C23").Value = 9.61471183610475
$ws.Range("D23").Value = 5.461355005830836
$ws.Range("F23").Value = 71.04972239975356
$ws.Range("G23").Value = 3.794684144046796
$ws.Range("J23").Value = 10.45166096748518
$ws.Range("L23").Value = 11.20681777767963
$ws.Range("M23").Value = 20.86134090705263
$ws.Range("B24").Value = 26.63683385303719
$ws.Range("C24").Value = 9.184702069197327
$ws.Range("D24").Value = 5.199419829738599
$ws.Range("F24").Value = 69.63362322170194
$ws.Range("G24").Value = 3.801992442415826
$ws.Range("J24").Value = 10.43830013652401
$ws.Range("L24").Value = 11.2624822379247
$ws.Range("M24").Value = 20.89056396696354
$ws.Range("B25").Value = 26.39963299128502
$ws.Range("C25").Value = 8.716589762319558
$ws.Range("D25").Value = 4.9019378594465
$ws.Range("F25").Value = 68.10368316137379
$ws.Range("G25").Value = 3.810422413716425
$ws.Range("J25").Value = 10.42375816753531
$ws.Range("L25").Value = 11.32890123517879
$ws.Range("M25").Value = 20.94229748936642
